$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$TextValue
    )
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $TextValue
    $cell.ClearFormats()
}

Set-TextValue 'D2' '64.286.38'
Set-TextValue 'E2' '  +0.84%  '
Set-TextValue 'D3' '3.501.02'
Set-TextValue 'E3' '  +0.50%  '
Set-TextValue 'E4' '  +0.15%  '
Set-TextValue 'D5' '586.12'
Set-TextValue 'E5' '  +0.70%  '
Set-TextValue 'D6' '134.20'
Set-TextValue 'E6' '  +3.30%  '
Set-TextValue 'E7' '  +0.02%  '
Set-TextValue 'E8' '  +1.05%  '
Set-TextValue 'E9' '  +1.40%  '
Set-TextValue 'E10' '  +2.21%  '
Set-TextValue 'E11' '  +2.40%  '
Set-TextValue 'D12' '4.099.34'
Set-TextValue 'E12' '  +1.09%  '
Set-TextValue 'E13' '  +3.87%  '
Set-TextValue 'E14' '  +1.43%  '
Set-TextValue 'D15' '3.501.06'
Set-TextValue 'E15' '  +0.49%  '
Set-TextValue 'E16' '  -4.09%  '
Set-TextValue 'D17' '64.286.43'
Set-TextValue 'E17' '  +0.92%  '
Set-TextValue 'D18' '9.89'
Set-TextValue 'E18' '  +0.73%  '
Set-TextValue 'E19' '  +2.81%  '
Set-TextValue 'D20' '13.62'
Set-TextValue 'E20' '  -2.88%  '
Set-TextValue 'D21' '393.52'
Set-TextValue 'E21' '  +3.73%  '
Set-TextValue 'E22' '  +0.07%  '
Set-TextValue 'D23' '3.641.56'
Set-TextValue 'E23' '  +0.61%  '
Set-TextValue 'D24' '74.28'
Set-TextValue 'E24' '  +1.67%  '
Set-TextValue 'E25' '  -0.03%  '
Set-TextValue 'E26' '  +2.35%  '
Set-TextValue 'D27' '1.00'
Set-TextValue 'E28' '  -0.36%  '
Set-TextValue 'B29' 'Fetch.AI'
Set-TextValue 'C29' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D29' '1.50'
Set-TextValue 'E29' '  -3.74%  '
Set-TextValue 'B30' 'InternetComputer(DFINITY)'
Set-TextValue 'C30' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D30' '8.29'
Set-TextValue 'E30' '  +1.27%  '
Set-TextValue 'E31' '  +1.41%  '
Set-TextValue 'D32' '3.520.75'
Set-TextValue 'E32' '  +0.92%  '
Set-TextValue 'E33' '  +4.81%  '
Set-TextValue 'E34' '  +0.05%  '
Set-TextValue 'E35' '  +0.63%  '
Set-TextValue 'D36' '5.15'
Set-TextValue 'E36' '  -1.99%  '
Set-TextValue 'E37' '  +1.14%  '
Set-TextValue 'D38' '6.91'
Set-TextValue 'E38' '  +0.32%  '
Set-TextValue 'D39' '164.30'
Set-TextValue 'E39' '  +2.75%  '
Set-TextValue 'D40' '0.0784'
Set-TextValue 'E40' '  -0.56%  '
Set-TextValue 'E41' '  -0.17%  '
Set-TextValue 'E42' '  +0.20%  '
Set-TextValue 'D43' '25.15'
Set-TextValue 'E43' '  -2.89%  '
Set-TextValue 'E44' '  +1.48%  '
Set-TextValue 'B45' 'ONDO'
Set-TextValue 'C45' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D45' '1.19'
Set-TextValue 'E45' '  -1.00%  '
Set-TextValue 'B46' 'Stacks'
Set-TextValue 'C46' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D46' '1.65'
Set-TextValue 'E46' '  +3.34%  '
Set-TextValue 'D47' '2.464.14'
Set-TextValue 'E47' '  +1.94%  '
Set-TextValue 'D48' '6.77'
Set-TextValue 'E48' '  -0.38%  '
Set-TextValue 'D49' '0.897'
Set-TextValue 'E49' '  +1.43%  '
Set-TextValue 'D50' '0.0261'
Set-TextValue 'E50' '  -0.48%  '
Set-TextValue 'E51' '  -0.02%  '

Write-Output "All cell updates applied."
